$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "Final report title: (the topic of your research.)"
#    -> "Final report title: Correlation Between Applicant Income And Approved Loan Amount"
#    (drop the parenthetical placeholder wording + its italics, keep a
#    normal, non-italic run with the real title)
# ---------------------------------------------------------------
$findTitle = $d.Content.Find
$findTitle.ClearFormatting()
$findTitle.Execute("the topic of your research.)", $false, $false, $false, $false, $false, $true, 1, $false, `
    "Correlation Between Applicant Income And Approved Loan Amount", 2) | Out-Null

$findTitleItalic = $d.Content.Find
$findTitleItalic.ClearFormatting()
if ($findTitleItalic.Execute("Correlation Between Applicant Income And Approved Loan Amount", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $findTitleItalic.Parent.Font.Italic = $false
}

$findParen = $d.Content.Find
$findParen.ClearFormatting()
$findParen.Execute("Final report title: (", $false, $false, $false, $false, $false, $true, 1, $false, `
    "Final report title: ", 2) | Out-Null

# ---------------------------------------------------------------
# 2) "Group ID:" -> "Group ID: A5"
# ---------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Group ID:") {
        $r = $p.Range
        $insertPoint = $d.Range($r.End - 1, $r.End - 1)
        $insertPoint.InsertAfter(" A5")
        break
    }
}

# ---------------------------------------------------------------
# 3) "Dataset number: " -> "Dataset number: DS256-loan_Data"
# ---------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Dataset number: ") {
        $r = $p.Range
        $insertPoint = $d.Range($r.End - 1, $r.End - 1)
        $insertPoint.InsertAfter("DS256-loan_Data")
        break
    }
}

# ---------------------------------------------------------------
# 4) Group-changes paragraph: merge the two runs that were split by the
#    stale "_GoBack" bookmark back into a single run (bookmark removed).
# ---------------------------------------------------------------
$groupText = "There were no changes to our group since the original allocation. " + `
    "We formed the group at the beginning of the project, and all members continued to work together throughout the entire process. " + `
    "Everyone remained committed, contributed consistently, and participated fully in the tasks assigned. " + `
    "Since no one left or joined the group, our GitHub access and member IDs stayed the same, and the group structure remained stable, allowing us to work smoothly without any disruptions"

$findGroup = $d.Content.Find
$findGroup.ClearFormatting()
if ($findGroup.Execute($groupText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $startPos = $findGroup.Parent.Start
    $endPos = $findGroup.Parent.End
    $rng = $d.Range($startPos, $endPos)
    $rng.Text = $groupText + "X"
    $fixEnd = $startPos + $groupText.Length + 1
    $rngFix = $d.Range($fixEnd - 1, $fixEnd)
    $rngFix.Text = ""
}

Write-Host "done"
